$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 13 (old "7455355 - Robson da Silva Rocha" row with no label),
# shifting rows 14-22 up to become rows 13-21.
$ws.Rows("13").Delete()

# Replace the long "Objetivos:" body text with the short placeholder text.
$ws.Range("B10").Value = "7455355 - Robson da Silva Rocha"
$ws.Range("C10").Value = "7455355 - Robson da Silva Rocha"

# Replace the long "Programa resumido:" body text (now row 13) with "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Replace the long "Programa:" body text (now row 15) with the date placeholder
# "01/01/2022". Copy/PasteSpecial(values) from the existing text cell B8 avoids
# Excel's automatic text->date conversion that a plain .Value assignment triggers.
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4163) | Out-Null
$ws.Range("C15").PasteSpecial(-4163) | Out-Null
$excel.Application.CutCopyMode = $false

# Replace the "Método:" body text (now row 18) with the short placeholder text.
$ws.Range("B18").Value = "7455355 - Robson da Silva Rocha"
$ws.Range("C18").Value = "7455355 - Robson da Silva Rocha"
